# Add new rows 201-236 to the Worksheet sheet for Alchemy of Souls
# season 1, episodes 15-16 quotes, matching the upstream commit
# 'Add files via upload'. Cell writes are ordered so that the
# generated shared-strings table lines up with the target workbook
# (Excel assigns shared-string indices in first-use order, and the
# original author entered a couple of cells out of strict row order).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: enter quote text (and the two brand-new character names)
# in the precise order the original strings were first typed, so new
# shared-string entries land at the same indices as the target file.
$ws.Cells.Item(202, 1).Value = 'The son of Cheonbugwan''s Gwanju should be protecting the royal family, but you brought in a wily beast instead.'
$ws.Cells.Item(201, 1).Value = 'Does this mean there is a soul shifter inside my palace?'
$ws.Cells.Item(203, 1).Value = 'I was told that Master Seo Gyeong destroyed all the ice stones 200 years ago after the Great Chaos.'
$ws.Cells.Item(204, 1).Value = 'It is not easy to quell human desires.'
$ws.Cells.Item(205, 1).Value = 'A power strong enough to defy the laws of nature and bring chaos to the world has returned once again. Thus, someone strong enough to stop it should also appear. That way, the world will be able to maintain balance.'
$ws.Cells.Item(206, 1).Value = 'Naksu may have been an assassin, but she never hurt ordinary people.'
$ws.Cells.Item(207, 1).Value = 'Naksu''s body was burned, and they all know that it was done by those with power and authority. Calling them fools when they are actually scared is negligence.'
$ws.Cells.Item(208, 1).Value = 'Your embarrassment will only be temporary, but I could really get killed for this.'
$ws.Cells.Item(209, 1).Value = 'You did not need to bring gifts to apologize for the past.'
$ws.Cells.Item(210, 1).Value = 'To be acknowledged as the daughter of the Jin family, you must open the door to Jinyowon.'
$ws.Cells.Item(211, 1).Value = 'A secret that ensures that we say the same thing.'
$ws.Cells.Item(212, 1).Value = 'If something like this ever happens again, and you find yourself in a life or death situation, give up on love and choose your life. That is an order from your master.'
$ws.Cells.Item(213, 1).Value = 'What if I do not mind dying? Then can I hold onto it and keep loving you?'
$ws.Cells.Item(214, 1).Value = 'I told you that I am willing to risk it all and that we are in this together. I will keep going no matter what. So you should do the same and stick to your confession. Keep loving me to death.'
$ws.Cells.Item(215, 1).Value = 'Master Lee knew that you were a soul shifter from the beginning. Perhaps, he even knows you are Naksu.'
$ws.Cells.Item(216, 1).Value = 'I heard that you took Gwigu, a relic of Jinyowon, to the palace today.'
$ws.Cells.Item(217, 1).Value = 'So this dog was Gwigu. It must have known that I was a soul shifter from the beginning. But it never barked at me. Was it covering for me too?'
$ws.Cells.Item(218, 1).Value = 'The rumors about soul shifters have spread throughout Daeho Fortress. If he hears the rumors, he will return, fearing for Songrim and for his son.'
$ws.Cells.Item(219, 1).Value = 'He is merciless towards those who stray.'
$ws.Cells.Item(220, 1).Value = 'And if you must draw and point that sword at someone, point it at me before at anyone else. If you stab me, then I will die as the one who released the assassin.'
$ws.Cells.Item(221, 1).Value = 'The heart is the most important. If you are in love, appearances do not matter.'
$ws.Cells.Item(222, 1).Value = 'My sister was her pride and joy. Even though she was blind, she could see energy and handle relics. Even at a young age, her abilities were as powerful as my mother''s.'
$ws.Cells.Item(223, 1).Value = 'Well, I heard that the Jin family members are born with more divine powers than spell-related powers.'
$ws.Cells.Item(224, 1).Value = 'A really powerful priestess can control people''s souls and trap them to be helpless.'
$ws.Cells.Item(225, 1).Value = 'If the priestess traps their powers, even the most powerful mage is helpless.'
$ws.Cells.Item(226, 1).Value = 'A Jinyowon priestess is in a different league.'
$ws.Cells.Item(227, 1).Value = 'He suffered because he knew what he was doing was wrong, and that fact was hard for him to endure.'
$ws.Cells.Item(228, 1).Value = 'To punish oneself more harshly than others would. Is it some sort of a heightened state of morality? Gosh, you were already there at that young age? That is incredible.'
$ws.Cells.Item(229, 1).Value = 'Is that right? I must have romanticized the memory. I thought we did the honorable thing. How embarrassing…'
$ws.Cells.Item(230, 1).Value = 'You should take my side like before. Come on, Uk.'
$ws.Cells.Item(231, 1).Value = 'I saw you wandering around until late. You usually go to bed like clockwork. Are you concerned about what is happening at Seoho Fortress?'
$ws.Cells.Item(232, 1).Value = 'You are feigning ignorance, right? It must be tough for you.'
$ws.Cells.Item(233, 1).Value = 'Each member of the Seo Family has a Lantern of Life. Because they guard the Seoho Fortress near the border and go to battle often, they each have a lantern that is lit by the energy of their soul so they know who is alive and dead.'
$ws.Cells.Item(233, 2).Value = 'Heo Yun-Ok'
$ws.Cells.Item(234, 2).Value = 'Shaman Bong'
$ws.Cells.Item(234, 1).Value = 'It is my business policy never to breach confidentiality. You can rest assured.'
$ws.Cells.Item(235, 1).Value = 'But there was a rumor that he married a woman from a powerful family.'
$ws.Cells.Item(236, 1).Value = 'Since you have Jin Cho-yeon''s blood in your body, you should be able to move the doors of Jinyowon. Try to move the doors even a little bit. That should be enough.'

# Step 2: fill in the remaining Character / Season / Episode cells.
$remaining = @(
    @{ Row = 201; Character = 'King Go Soon'; Season = 1; Episode = 15 }
    @{ Row = 202; Character = 'Jin Mu'; Season = 1; Episode = 15 }
    @{ Row = 203; Character = 'Jin Mu'; Season = 1; Episode = 15 }
    @{ Row = 204; Character = 'Master Lee'; Season = 1; Episode = 15 }
    @{ Row = 205; Character = 'Master Lee'; Season = 1; Episode = 15 }
    @{ Row = 206; Character = 'Seo Yul'; Season = 1; Episode = 15 }
    @{ Row = 207; Character = 'Jang Uk'; Season = 1; Episode = 15 }
    @{ Row = 208; Character = 'Mu-deok / Naksu'; Season = 1; Episode = 15 }
    @{ Row = 209; Character = 'Lady Jin'; Season = 1; Episode = 15 }
    @{ Row = 210; Character = 'Jin Mu'; Season = 1; Episode = 15 }
    @{ Row = 211; Character = 'Jang Uk'; Season = 1; Episode = 15 }
    @{ Row = 212; Character = 'Mu-deok / Naksu'; Season = 1; Episode = 15 }
    @{ Row = 213; Character = 'Jang Uk'; Season = 1; Episode = 15 }
    @{ Row = 214; Character = 'Jang Uk'; Season = 1; Episode = 15 }
    @{ Row = 215; Character = 'Jang Uk'; Season = 1; Episode = 16 }
    @{ Row = 216; Character = 'Lady Jin'; Season = 1; Episode = 16 }
    @{ Row = 217; Character = 'Mu-deok / Naksu'; Season = 1; Episode = 16 }
    @{ Row = 218; Character = 'Park Jin'; Season = 1; Episode = 16 }
    @{ Row = 219; Character = 'Heo Yeom'; Season = 1; Episode = 16 }
    @{ Row = 220; Character = 'Jang Uk'; Season = 1; Episode = 16 }
    @{ Row = 221; Character = 'Kim Do-Joo'; Season = 1; Episode = 16 }
    @{ Row = 222; Character = 'Jin Cho-Yeon'; Season = 1; Episode = 16 }
    @{ Row = 223; Character = 'Mu-deok / Naksu'; Season = 1; Episode = 16 }
    @{ Row = 224; Character = 'Jin Cho-Yeon'; Season = 1; Episode = 16 }
    @{ Row = 225; Character = 'Jin Cho-Yeon'; Season = 1; Episode = 16 }
    @{ Row = 226; Character = 'Jin Cho-Yeon'; Season = 1; Episode = 16 }
    @{ Row = 227; Character = 'Jang Uk'; Season = 1; Episode = 16 }
    @{ Row = 228; Character = 'Park Dang-Gu'; Season = 1; Episode = 16 }
    @{ Row = 229; Character = 'Park Dang-Gu'; Season = 1; Episode = 16 }
    @{ Row = 230; Character = 'Park Dang-Gu'; Season = 1; Episode = 16 }
    @{ Row = 231; Character = 'Park Dang-Gu'; Season = 1; Episode = 16 }
    @{ Row = 232; Character = 'Jang Uk'; Season = 1; Episode = 16 }
    @{ Row = 233; Character = 'Heo Yun-Ok'; Season = 1; Episode = 16 }
    @{ Row = 234; Character = 'Shaman Bong'; Season = 1; Episode = 16 }
    @{ Row = 235; Character = 'Shaman Bong'; Season = 1; Episode = 16 }
    @{ Row = 236; Character = 'Jin Mu'; Season = 1; Episode = 16 }
)
foreach ($r in $remaining) {
    $ws.Cells.Item($r.Row, 2).Value = $r.Character
    $ws.Cells.Item($r.Row, 3).Value = $r.Season
    $ws.Cells.Item($r.Row, 4).Value = $r.Episode
}

# Match the final cell selection recorded in the saved workbook
$ws.Range("A251").Select()

